# Auto-generated edit script applying numeric corrections
# described by the upstream diff for Adamantoise_Profits (Leve profit sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1295.4559
$ws.Range("I15").Value = 1295.4559
$ws.Range("K15").Value = 3886.3677
$ws.Range("M15").Value = -3717.3677
$ws.Range("H33").Value = 378.4516
$ws.Range("I33").Value = 227.04347
$ws.Range("J33").Value = 813.75
$ws.Range("K33").Value = 227.04347
$ws.Range("L33").Value = 813.75
$ws.Range("M33").Value = 1.956529999999987
$ws.Range("N33").Value = -1271.75
$ws.Range("H62").Value = 8059
$ws.Range("J62").Value = 9756
$ws.Range("L62").Value = 9756
$ws.Range("N62").Value = -11004
$ws.Range("H65").Value = 8059
$ws.Range("J65").Value = 9756
$ws.Range("L65").Value = 48780
$ws.Range("N65").Value = -55020
$ws.Range("H137").Value = 4059.1614
$ws.Range("I137").Value = 3518.3333
$ws.Range("K137").Value = 10554.9999
$ws.Range("M137").Value = -8004.999899999999
$ws.Range("H141").Value = 5041.3335
$ws.Range("I141").Value = 4330
$ws.Range("K141").Value = 12990
$ws.Range("M141").Value = -7810

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13024193
$ws.Range("I32").Value = 7937818.5
$ws.Range("K32").Value = 7937818.5
$ws.Range("M32").Value = -7937531.5
$ws.Range("H102").Value = 2928.1667
$ws.Range("J102").Value = 2916.6667
$ws.Range("L102").Value = 2916.6667
$ws.Range("N102").Value = -6160.6667
$ws.Range("H132").Value = 2850.4565
$ws.Range("I132").Value = 2473.5881
$ws.Range("K132").Value = 7420.7643
$ws.Range("M132").Value = -4890.7643

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 105313.5
$ws.Range("J60").Value = 105313.5
$ws.Range("L60").Value = 105313.5
$ws.Range("N60").Value = -106511.5
$ws.Range("H94").Value = 1328.75
$ws.Range("I94").Value = 647.3077
$ws.Range("K94").Value = 647.3077
$ws.Range("M94").Value = -196.3077

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3628.975
$ws.Range("I31").Value = 2046.0588
$ws.Range("J31").Value = 4798.9565
$ws.Range("K31").Value = 2046.0588
$ws.Range("L31").Value = 4798.9565
$ws.Range("M31").Value = -1751.0588
$ws.Range("N31").Value = -5388.9565
$ws.Range("H34").Value = 3628.975
$ws.Range("I34").Value = 2046.0588
$ws.Range("J34").Value = 4798.9565
$ws.Range("K34").Value = 2046.0588
$ws.Range("L34").Value = 4798.9565
$ws.Range("M34").Value = -1844.0588
$ws.Range("N34").Value = -5202.9565
$ws.Range("H64").Value = 107000
$ws.Range("J64").Value = 107000
$ws.Range("L64").Value = 107000
$ws.Range("N64").Value = -107496
$ws.Range("H67").Value = 107000
$ws.Range("J67").Value = 107000
$ws.Range("L67").Value = 107000
$ws.Range("N67").Value = -108716
$ws.Range("H68").Value = 123000
$ws.Range("J68").Value = 123000
$ws.Range("L68").Value = 123000
$ws.Range("N68").Value = -124498
$ws.Range("H71").Value = 123000
$ws.Range("J71").Value = 123000
$ws.Range("L71").Value = 369000
$ws.Range("N71").Value = -376488
$ws.Range("H82").Value = 40000
$ws.Range("J82").Value = 40000
$ws.Range("L82").Value = 40000
$ws.Range("N82").Value = -40722
$ws.Range("H85").Value = 40000
$ws.Range("J85").Value = 40000
$ws.Range("L85").Value = 40000
$ws.Range("N85").Value = -42496
$ws.Range("H105").Value = 2429.75
$ws.Range("I105").Value = 2439.6667
$ws.Range("J105").Value = 2400
$ws.Range("K105").Value = 2439.6667
$ws.Range("L105").Value = 2400
$ws.Range("M105").Value = -692.6667000000002
$ws.Range("N105").Value = -5894
$ws.Range("H110").Value = 82319.664
$ws.Range("J110").Value = 82319.664
$ws.Range("L110").Value = 82319.664
$ws.Range("N110").Value = -90499.664
$ws.Range("H122").Value = 2961.4614
$ws.Range("I122").Value = 2633.111
$ws.Range("K122").Value = 7899.333
$ws.Range("M122").Value = -5449.333
$ws.Range("H132").Value = 2152
$ws.Range("I132").Value = 1713.8462
$ws.Range("K132").Value = 5141.5386
$ws.Range("M132").Value = -2611.5386
$ws.Range("H134").Value = 2282.5
$ws.Range("I134").Value = 1250.875
$ws.Range("J134").Value = 4758.4
$ws.Range("K134").Value = 3752.625
$ws.Range("L134").Value = 14275.2
$ws.Range("M134").Value = -1217.625
$ws.Range("N134").Value = -19345.2
$ws.Range("H141").Value = 523086
$ws.Range("J141").Value = 568062.1
$ws.Range("L141").Value = 568062.1
$ws.Range("N141").Value = -578422.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 183.6
$ws.Range("I86").Value = 81.8
$ws.Range("J86").Value = 285.4
$ws.Range("K86").Value = 245.4
$ws.Range("L86").Value = 856.1999999999999
$ws.Range("M86").Value = 940.6
$ws.Range("N86").Value = -3228.2
$ws.Range("H89").Value = 183.6
$ws.Range("I89").Value = 81.8
$ws.Range("J89").Value = 285.4
$ws.Range("K89").Value = 736.1999999999999
$ws.Range("L89").Value = 2568.6
$ws.Range("M89").Value = 5191.8
$ws.Range("N89").Value = -14424.6
$ws.Range("H107").Value = 586.5
$ws.Range("J107").Value = 609.3570999999999
$ws.Range("L107").Value = 1828.0713
$ws.Range("N107").Value = -5668.0713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2400.25
$ws.Range("I113").Value = 2400.25
$ws.Range("K113").Value = 2400.25
$ws.Range("M113").Value = -230.25
$ws.Range("H132").Value = 3074.75
$ws.Range("I132").Value = 2899.7273
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 8699.1819
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -6169.1819
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 8900
$ws.Range("J14").Value = 8900
$ws.Range("L14").Value = 8900
$ws.Range("N14").Value = -9244
$ws.Range("H26").Value = 6502.25
$ws.Range("I26").Value = 6502.25
$ws.Range("K26").Value = 6502.25
$ws.Range("M26").Value = -6207.25
$ws.Range("H31").Value = 3140.9
$ws.Range("I31").Value = 1759
$ws.Range("J31").Value = 3601.5334
$ws.Range("K31").Value = 1759
$ws.Range("L31").Value = 3601.5334
$ws.Range("M31").Value = -1511
$ws.Range("N31").Value = -4097.5334
$ws.Range("H40").Value = 8880.333000000001
$ws.Range("I40").Value = 8880.333000000001
$ws.Range("K40").Value = 8880.333000000001
$ws.Range("M40").Value = -8744.333000000001
$ws.Range("H46").Value = 3315.7
$ws.Range("I46").Value = 1949.5
$ws.Range("K46").Value = 1949.5
$ws.Range("M46").Value = -1761.5
$ws.Range("H56").Value = 30025.5
$ws.Range("I56").Value = 23700.666
$ws.Range("K56").Value = 23700.666
$ws.Range("M56").Value = -23009.666
$ws.Range("H122").Value = 2931.76
$ws.Range("I122").Value = 3002.1365
$ws.Range("J122").Value = 2415.6667
$ws.Range("K122").Value = 9006.4095
$ws.Range("L122").Value = 7247.000100000001
$ws.Range("M122").Value = -6556.4095
$ws.Range("N122").Value = -12147.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 14901.444
$ws.Range("J45").Value = 18230.428
$ws.Range("L45").Value = 18230.428
$ws.Range("N45").Value = -19212.428
$ws.Range("H59").Value = 100102
$ws.Range("J59").Value = 100102
$ws.Range("L59").Value = 100102
$ws.Range("N59").Value = -101578
